$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The body placeholder shape on slide 1 (Google Shape;69;p13) holding the
# bulleted "remainder of today's class" instructions.
$shp = $s.Shapes.Item(2)

# Widen the text box slightly (cx: 4128000 -> 4260300 EMU).
$shp.Width = 4260300 / 12700

$tf = $shp.TextFrame
$tr = $tf.TextRange

# Replace the final bullet's text in place (keeps its run/paragraph formatting).
$oldText = "If there is any time remaining, you will begin work on your homework assignment."
$newText = "Write your name on each completed sheet."
$full = $tr.Text
$idx = $full.IndexOf($oldText)
if ($idx -ge 0) {
    $target = $tr.Characters($idx + 1, $oldText.Length)
    $target.Text = $newText
}

# Append a brand-new bullet after it, inheriting the same paragraph formatting.
$tr2 = $tf.TextRange
$null = $tr2.InsertAfter([char]13 + "Submit to the designated MyCourses" + [char]0x2019 + " Dropbox.")
